$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$table = $ws.ListObjects.Item("Tabelle1")

# Add three new rows to the table so that it (and its underlying range) grows
# from A1:I8 to A1:I11.
$table.ListRows.Add() | Out-Null
$table.ListRows.Add() | Out-Null
$table.ListRows.Add() | Out-Null

# Column A - File Name
$ws.Range("A9").Value = "2023-02-16-2049_log_reg_centers.csv"
$ws.Range("A10").Value = "2023-02-16-2136_xgb_centers.csv"
$ws.Range("A11").Value = "2023-02-21-2055_xgb_centers_nohyp.csv"

# Column B - Model
$ws.Range("B9").Value = "Logistic Classifier"
$ws.Range("B10").Value = "XGB Classifier"
$ws.Range("B11").Value = "XGB Classifier"

# Column C - Data
$ws.Range("C9").Value = "MoCov"
$ws.Range("C10").Value = "MoCov"
$ws.Range("C11").Value = "MoCov"

# Column D - Parameters
$ws.Range("D9").Value = "-"
$ws.Range("D10").Value = "{colsample_bylevel=0.2}"
$ws.Range("D11").Value = "-"

# Column E - repeated CV
$ws.Range("E9").Value = "1 x 3"
$ws.Range("E10").Value = "1 x 3"
$ws.Range("E11").Value = "1 x 3"

# Column F - Grade Weakly Supervision
$ws.Range("F9").Value = "weakly supervision with cv centers"
$ws.Range("F10").Value = "weakly supervision with cv centers"
$ws.Range("F11").Value = "weakly supervision with cv centers"

# Column H - Hand in
$ws.Range("H9").Value = "Feb. 16, 2023, 7:51 p.m."
$ws.Range("H10").Value = "Feb. 16, 2023, 8:37 p.m."
$ws.Range("H11").Value = "Feb. 21, 2023, 7:58 p.m."

# Column G - Average AUC (entered bottom-up)
$ws.Range("G11").Value = "0.601 (0.019)"
$ws.Range("G10").Value = "0.597 (0.003)"
$ws.Range("G9").Value = "0.608 (0.029)"

# Column I - AUC Test
$ws.Range("I9").Value = 0.631
$ws.Range("I10").Value = 0.587
$ws.Range("I11").Value = 0.599

$ws.Range("G10").Select()
